$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final keyword/appID table for rows 2-18 (row 1 header unchanged):
# row -> (rowHeight-or-$null-if-unchanged, keyword, appID)
$rows = @(
    @(2,  $null, "stretchy taxi",        "com.singleton.strechy"),
    @(3,  24,    "passive income",       "passive.income.nadi.myfirstdrawermenuproject"),
    @(4,  12.8,  "taxi game",            "com.singleton.strechy"),
    @(5,  12.8,  "passive income",       "passive.income.nadi.myfirstdrawermenuproject"),
    @(6,  $null, "bitcoin",              "com.hamxa.shaynachim"),
    @(7,  24,    "passive income ideas", "passive.income.nadi.myfirstdrawermenuproject"),
    @(8,  $null, "taxi free game",       "com.singleton.strechy"),
    @(9,  12.8,  "best bitcoin",         "com.hamxa.shaynachim"),
    @(10, $null, "bitcoin course",       "com.hamxa.shaynachim"),
    @(11, 12.8,  "taxi",                 "com.singleton.strechy"),
    @(12, 24,    "bitcoin for beginners  ", "com.hamxa.shaynachim"),
    @(13, 24,    "bitcoin beginners  ",  "com.hamxa.shaynachim"),
    @(14, $null, "taxi game free",       "com.singleton.strechy"),
    @(15, 24,    "taxi free challenge",  "com.singleton.strechy"),
    @(16, $null, "taxi offline game",    "com.singleton.strechy"),
    @(17, $null, "challenge taxi game",  "com.singleton.strechy")
)

foreach ($row in $rows) {
    $r = $row[0]
    $h = $row[1]
    $kw = $row[2]
    $appid = $row[3]
    $ws.Cells.Item($r, 1).Value = $kw
    $ws.Cells.Item($r, 2).Value = $appid
    if ($h -ne $null) {
        $ws.Rows.Item($r).RowHeight = $h
    }
}

# Add new row 18 (copy formatting from row 17, then set values + height)
$ws.Range("A17:B17").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)
$ws.Cells.Item(18, 1).Value = "offline taxi game"
$ws.Cells.Item(18, 2).Value = "com.singleton.strechy"
$ws.Rows.Item(18).RowHeight = 24

# Update selection to match the post-edit state (A17:A18 selected)
$ws.Range("A17:A18").Select()
